# unify the conception of DataNode, DataTable, Entity.
# Rename the worksheet that used to model a loose "Property" bag so it
# matches the unified DataNode/DataTable/Entity naming used elsewhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core semantic edit: rename the sheet -------------------------------
$ws.Name = "DataNode"

# --- Cosmetic follow-up: the editor's cursor ended up on D40 -----------
$ws.Range("D40").Select()

# --- Minor column-width touch-ups (re-saved file nudged a few widths) --
$ws.Columns.Item(1).ColumnWidth = 18.428571428571427   # -> ~19.125 "A"
$ws.Columns.Item(2).ColumnWidth = 7.428571428571429    # -> ~8.125  "B"
$ws.Columns.Item(3).ColumnWidth = 27.428571428571427   # -> ~28.125 "C"
$ws.Columns.Item(6).ColumnWidth = 11.857142857142858   # -> ~12.625 "F"
$ws.Columns.Item(7).ColumnWidth = 11.857142857142858   # -> ~12.625 "G"

# --- Add the small 9pt phonetic-guide font and flag the sheet as using
#     no IME phonetic conversion (matches the Windows-build resave). ----
$ws.Range("A1").Phonetics.Font.Size = 9
$ws.Range("A1").Phonetics.Font.Name = "宋体"
